# Updated cryptos list on Sat Jun 17 21:46:56 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    # Force the value to be stored as text (avoids Excel auto-converting
    # number-looking strings like "0.9997" to a numeric type), while
    # restoring the cell's original style so no quote-prefix styling sticks.
    $origStyle = $cell.Style
    $cell.Value = "'" + $text
    $cell.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "26.504.07"
Set-TextValue $ws.Range("E2") "  +0.74%  "
Set-TextValue $ws.Range("D3") "1.727.15"
Set-TextValue $ws.Range("E3") "  +0.52%  "
Set-TextValue $ws.Range("D4") "0.9997"
Set-TextValue $ws.Range("E4") "  -0.10%  "
Set-TextValue $ws.Range("D5") "245.18"
Set-TextValue $ws.Range("E5") "  +2.69%  "
Set-TextValue $ws.Range("E6") "  -0.11%  "
Set-TextValue $ws.Range("D7") "0.4799"
Set-TextValue $ws.Range("E7") "  +2.03%  "
Set-TextValue $ws.Range("D8") "0.2661"
Set-TextValue $ws.Range("E8") "  +1.44%  "
Set-TextValue $ws.Range("D9") "0.06214"
Set-TextValue $ws.Range("E9") "  +0.48%  "
Set-TextValue $ws.Range("D10") "1.725.21"
Set-TextValue $ws.Range("E10") "  +0.44%  "
Set-TextValue $ws.Range("D11") "0.07147"
Set-TextValue $ws.Range("E11") "  +1.12%  "
Set-TextValue $ws.Range("E12") "  +2.39%  "
Set-TextValue $ws.Range("D13") "0.6154"
Set-TextValue $ws.Range("E13") "  +4.19%  "
Set-TextValue $ws.Range("D14") "4.513"
Set-TextValue $ws.Range("E14") "  +2.89%  "
Set-TextValue $ws.Range("D15") "77.09"
Set-TextValue $ws.Range("E15") "  +1.34%  "
Set-TextValue $ws.Range("E16") "  -0.08%  "
Set-TextValue $ws.Range("D17") "26.503.62"
Set-TextValue $ws.Range("E17") "  +0.73%  "
Set-TextValue $ws.Range("E18") "  -0.15%  "
Set-TextValue $ws.Range("D19") "0.000006929"
Set-TextValue $ws.Range("E19") "  +2.01%  "
Set-TextValue $ws.Range("E20") "  +0.88%  "
Set-TextValue $ws.Range("D21") "1.946.79"
Set-TextValue $ws.Range("E21") "  +0.37%  "
Set-TextValue $ws.Range("D22") "4.514"
Set-TextValue $ws.Range("E22") "  -0.57%  "
Set-TextValue $ws.Range("D23") "8.938"
Set-TextValue $ws.Range("E23") "  +2.37%  "
Set-TextValue $ws.Range("D24") "5.273"
Set-TextValue $ws.Range("E24") "  -0.83%  "
Set-TextValue $ws.Range("E25") "  +0.72%  "
Set-TextValue $ws.Range("D26") "15.33"
Set-TextValue $ws.Range("E26") "  +0.64%  "
Set-TextValue $ws.Range("D27") "1.793"
Set-TextValue $ws.Range("E27") "  +2.03%  "
Set-TextValue $ws.Range("E28") "  +0.22%  "
Set-TextValue $ws.Range("D29") "106.83"
Set-TextValue $ws.Range("E29") "  -1.19%  "
Set-TextValue $ws.Range("D30") "3.971"
Set-TextValue $ws.Range("E30") "  -0.74%  "
Set-TextValue $ws.Range("D31") "0.08023"
Set-TextValue $ws.Range("E31") "  +3.90%  "
Set-TextValue $ws.Range("D32") "3.705"
Set-TextValue $ws.Range("E32") "  +0.72%  "
Set-TextValue $ws.Range("D33") "0.04562"
Set-TextValue $ws.Range("E33") "  +2.66%  "
Set-TextValue $ws.Range("D34") "0.9997"
Set-TextValue $ws.Range("E34") "  -0.11%  "
Set-TextValue $ws.Range("E35") "  +0.00%  "
Set-TextValue $ws.Range("D36") "0.6336"
Set-TextValue $ws.Range("E36") "  +2.40%  "
Set-TextValue $ws.Range("D37") "0.9910"
Set-TextValue $ws.Range("E37") "  +1.68%  "
Set-TextValue $ws.Range("D38") "0.9251"
Set-TextValue $ws.Range("E38") "  +0.28%  "
Set-TextValue $ws.Range("D39") "2.089"
Set-TextValue $ws.Range("E39") "  +10.27%  "
Set-TextValue $ws.Range("E40") "  +0.46%  "
Set-TextValue $ws.Range("D41") "105.04"
Set-TextValue $ws.Range("E41") "  -7.97%  "
Set-TextValue $ws.Range("D42") "1.005"
Set-TextValue $ws.Range("E42") "  +0.38%  "
Set-TextValue $ws.Range("D43") "0.01499"
Set-TextValue $ws.Range("E43") "  +1.24%  "
Set-TextValue $ws.Range("D44") "5.571"
Set-TextValue $ws.Range("E44") "  +4.26%  "
Set-TextValue $ws.Range("D45") "0.3888"
Set-TextValue $ws.Range("E45") "  +2.18%  "
Set-TextValue $ws.Range("D46") "6.947"
Set-TextValue $ws.Range("E46") "  +10.96%  "
Set-TextValue $ws.Range("D47") "0.1181"
Set-TextValue $ws.Range("E47") "  +1.64%  "
Set-TextValue $ws.Range("D48") "0.05332"
Set-TextValue $ws.Range("E48") "  +0.79%  "
Set-TextValue $ws.Range("D49") "30.87"
Set-TextValue $ws.Range("E49") "  +1.21%  "
Set-TextValue $ws.Range("D50") "7.822"
Set-TextValue $ws.Range("E50") "  +1.53%  "
Set-TextValue $ws.Range("E51") "  +4.31%  "
